$d = $word.ActiveDocument

# 1. Mark the run that carries the inline picture as "do not spell/grammar
#    check" (w:rPr/w:noProof) -- this is what Word stamps on a run as soon
#    as it contains a picture that was inserted/pasted into the document.
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs($i)
    if ($para.Range.InlineShapes.Count -gt 0) {
        $para.Range.NoProofing = $true
    }
}

# 2. Add six more blank paragraphs right after the picture paragraph (i.e.
#    just before the already-existing trailing blank paragraph), mirroring
#    a student pressing Enter a few extra times while finishing up the lab.
$wNs = "http://schemas.openxmlformats.org/wordprocessingml/2006/main"
for ($i = 0; $i -lt 6; $i++) {
    $lastPara = $d.Paragraphs.Last
    $insertionPoint = $d.Range($lastPara.Range.Start, $lastPara.Range.Start)
    $insertionPoint.InsertXML("<w:p xmlns:w='$wNs'/>") | Out-Null
}
